# fix: data model and list
#
# 1) Row 9 ("Anubis"): add a Description, update Specie List from
#    "European Shorthair" to "Bengal Cat", and set Gender List to "Male".
# 2) Rows 11-21: fill in the Weight (column J) with 0.05 for every row.
# 3) Move the visible selection to J10:J21 (matches the author's last
#    selection after filling the Weight column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: Anubis -------------------------------------------------------
$ws.Range("D9").Value = "This male Bengal cat is called Anubis. He is very curious, but shy."
$ws.Range("D9").WrapText = $true

$ws.Range("K9").Value = "Bengal Cat"

$ws.Range("L9").Value = "Male"

# --- Rows 11-21: Weight column (J) = 0.05 --------------------------------
for ($r = 11; $r -le 21; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.05
}

# --- Final selection -------------------------------------------------------
$null = $ws.Range("J10:J21").Select()
